# Version finale Fiche de paie Anka avec l'ajout de Autre allocation
#
# The original payslip has a "TOTAL DES RETENUES" line at row 51. This
# edit inserts a new "Autre allocation" line above it (new row 51),
# pushing "TOTAL DES RETENUES" and everything below it down by one row.
# The new allocation (20000) is added to the net pay which is now stored
# as a plain value instead of the old rounding formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above the old row 51 ("TOTAL DES RETENUES") ---
# This shifts rows 51-58 down to 52-59 and automatically repoints the
# dependent formula in E49 (=+E43-D51 -> =+E43-D52).
$ws.Rows("51:51").Insert()

# Reproduce the formatting of the row above (row 50) on the new row 51,
# including the thin left/right borders used throughout this block and
# the 15.75pt row height used by all the other rows in the table.
$ws.Range("A50:I50").Copy()
$ws.Range("A51:I51").PasteSpecial(-4122)
$ws.Rows(51).RowHeight = 15.75

# --- 2. Fill in the new "Autre allocation" row ---
$ws.Range("B51").Value = "Autre allocation"
$ws.Range("C51").Value = ""
$ws.Range("D51").Value = ""
$ws.Range("E51").Value = 20000

# --- 3. The "NET A PAYER" cell (old row 53, now row 54) is converted
# from a rounding formula to the new total as a plain number
# (180794 + 20000 = 200794). ---
$ws.Range("E54").Value = 200794

# --- 4. Restore the active cell selection as saved in the workbook ---
[void]$ws.Range("I49").Select()
